# Update the cryptos worksheet with refreshed price/volume data (and a
# re-ordering of a few coin rows) per the latest symbol-list scrape.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Text($cellRef, $text) {
    # Leading apostrophe forces Excel to store the value as literal text
    # (preserving things like trailing zeros / leading '-' / '%' signs)
    # instead of silently re-interpreting it as a number or percentage.
    $ws.Range($cellRef).Value = "'" + $text
}

# Row 2 - BNB
Set-Text "D2" "329.10"
Set-Text "E2" "6.46%"

# Row 3 - OKB
Set-Text "D3" "40.21"
Set-Text "E3" "7.59%"

# Row 4 - HuobiToken
Set-Text "D4" "5.565"
Set-Text "E4" "8.43%"

# Row 5 - Cronos
Set-Text "E5" "3.54%"

# Row 6 - GateToken
Set-Text "D6" "4.550"
Set-Text "E6" "3.43%"

# Row 7 - was FTXToken, now KuCoinToken
Set-Text "B7" "KuCoinToken"
Set-Text "C7" "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
Set-Text "D7" "8.676"
Set-Text "E7" "5.05%"

# Row 8 - was KuCoinToken, now FTXToken
Set-Text "B8" "FTXToken"
Set-Text "C8" "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-Text "D8" "1.990"
Set-Text "E8" "5.81%"

# Row 9 - BTSEToken
Set-Text "E9" "0.19%"

# Row 10 - MXToken
Set-Text "D10" "0.9491"
Set-Text "E10" "2.58%"

# Row 11 - LiechtensteinCryptoassetsExchange
Set-Text "D11" "0.1275"
Set-Text "E11" "14.52%"

# Row 12 - WazirX
Set-Text "D12" "0.1974"
Set-Text "E12" "3.23%"

# Row 13 - MandalaExchangeToken
Set-Text "D13" "0.09180"
Set-Text "E13" "3.36%"

# Row 14 - BitrueCoin
Set-Text "D14" "0.03585"
Set-Text "E14" "7.24%"

# Row 15 - BitMartToken
Set-Text "D15" "0.09596"
Set-Text "E15" "-0.11%"

# Row 16 - BitForexToken
Set-Text "D16" "0.001305"
Set-Text "E16" "-5.39%"

# Row 17 - TigerCash
Set-Text "D17" "0.006094"
Set-Text "E17" "1.45%"

# Row 18 - was HotbitToken, now LEO
Set-Text "B18" "LEO"
Set-Text "C18" "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-Text "D18" "3.366"
Set-Text "E18" "-0.85%"

# Row 19 - was LEO, now BitpandaEcosystemToken
Set-Text "B19" "BitpandaEcosystemToken"
Set-Text "C19" "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
Set-Text "D19" "0.3508"
Set-Text "E19" "1.50%"

# Row 20 - was BitpandaEcosystemToken, now MCDex
Set-Text "B20" "MCDex"
Set-Text "C20" "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-Text "D20" "7.465"
Set-Text "E20" "17.26%"

# Row 21 - was MCDex, now ProBitToken
Set-Text "B21" "ProBitToken"
Set-Text "C21" "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
Set-Text "D21" "0.1351"
Set-Text "E21" "2.81%"

# Row 22 - was ProBitToken, now ZBToken
Set-Text "B22" "ZBToken"
Set-Text "C22" "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
Set-Text "D22" "0.2490"
Set-Text "E22" "3.65%"

# Row 23 - was ZBToken, now CoinExToken
Set-Text "B23" "CoinExToken"
Set-Text "C23" "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-Text "D23" "0.04424"
Set-Text "E23" "1.70%"

# Row 24 - was CoinExToken, now BitKan
Set-Text "B24" "BitKan"
Set-Text "C24" "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
Set-Text "D24" "0.001227"
Set-Text "E24" "2.12%"

# Row 25 - was BitKan, now HotbitToken
Set-Text "B25" "HotbitToken"
Set-Text "C25" "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
Set-Text "D25" "0.004278"
Set-Text "E25" "-0.01%"

# Row 26 - NitroEx
Set-Text "E26" "-14.35%"

# Row 27 - UpBots
Set-Text "D27" "0.0003993"
Set-Text "E27" "37.50%"

# Row 39 - One
Set-Text "D39" "0.02514"
Set-Text "E39" "16.41%"

# Row 40 - IDEX
Set-Text "D40" "0.05222"
Set-Text "E40" "4.04%"

# Row 41 - KickToken
Set-Text "D41" "0.007818"
Set-Text "E41" "3.13%"

# Row 42 - BKEXToken
Set-Text "E42" "5.74%"

# Row 43 - Dexo
Set-Text "D43" "0.008894"
Set-Text "E43" "4.40%"

# Row 44 - CEJI
Set-Text "E44" "6.58%"

# Row 45 - LocalTraders
Set-Text "D45" "0.009611"
Set-Text "E45" "17.89%"

# Row 46 - CoinLion
Set-Text "D46" "0.00006664"
Set-Text "E46" "2.25%"

# Row 47 - Kangarootoken
Set-Text "E47" "-0.06%"

# Row 48 - BOLO
Set-Text "D48" "0.002904"
Set-Text "E48" "-11.95%"

# Row 49 - CoinbaseStockToken
Set-Text "D49" "0.002302"
Set-Text "E49" "59.31%"

# Row 50 - CryptobidCoin
Set-Text "E50" "-0.06%"

# Row 51 - SpecialPowerGold
Set-Text "D51" "0.0002001"
Set-Text "E51" "-0.06%"

Write-Output "cryptos sheet updated"
